# Insert a new data row (2023) above the existing "2025" row, shifting
# the existing rows down. After the shift, update the (now) 2025 row's
# probability values to match what used to be on the 2030 row, and give
# the new 2023 row the probability values that used to be on the 2025 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2:3 down to 3:4, freeing up row 2 for the new "2023" entry.
$ws.Rows.Item(2).Insert()
# Inserting copies formatting from the row above (the bold/bordered header
# row) - strip that back off so the new data row matches its siblings.
$ws.Range("A2:D2").ClearFormats()

# New row 2: 2023, with the probabilities the 2025 row used to carry.
$ws.Cells.Item(2, 1).Value = 2023
$ws.Cells.Item(2, 2).Value = 0.924
$ws.Cells.Item(2, 3).Value = 0.06
$ws.Cells.Item(2, 4).Value = 0.016

# Row 3 (previously row 2, "2025") now takes the probabilities that used
# to belong to the 2030 row.
$ws.Cells.Item(3, 1).Value = 2025
$ws.Cells.Item(3, 2).Value = 0.81
$ws.Cells.Item(3, 3).Value = 0.15
$ws.Cells.Item(3, 4).Value = 0.04

# Row 4 (previously row 3, "2030") keeps its original probabilities.
$ws.Cells.Item(4, 1).Value = 2030
$ws.Cells.Item(4, 2).Value = 0.81
$ws.Cells.Item(4, 3).Value = 0.15
$ws.Cells.Item(4, 4).Value = 0.04
